$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value = "Bitcoin"
$c.Style = "Normal"
$c = $ws.Range("C2")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$c.Style = "Normal"
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "51.366.71"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -1.78%  "
$c.Style = "Normal"

$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "Ethereum"
$c.Style = "Normal"
$c = $ws.Range("C3")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.926.86"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -2.30%  "
$c.Style = "Normal"

$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "TetherUSD"
$c.Style = "Normal"
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.21%  "
$c.Style = "Normal"

$c = $ws.Range("B5")
$c.NumberFormat = "@"
$c.Value = "BNB"
$c.Style = "Normal"
$c = $ws.Range("C5")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "374.33"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +5.63%  "
$c.Style = "Normal"

$c = $ws.Range("B6")
$c.NumberFormat = "@"
$c.Value = "Solana"
$c.Style = "Normal"
$c = $ws.Range("C6")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "104.09"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -4.09%  "
$c.Style = "Normal"

$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = "XRP"
$c.Style = "Normal"
$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.548"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -2.95%  "
$c.Style = "Normal"

$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = "USDC"
$c.Style = "Normal"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = "Cardano"
$c.Style = "Normal"
$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.593"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -5.03%  "
$c.Style = "Normal"

$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c.Style = "Normal"
$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "37.31"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -3.16%  "
$c.Style = "Normal"

$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = "TRON"
$c.Style = "Normal"
$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.139"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -0.31%  "
$c.Style = "Normal"

$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = "Dogecoin"
$c.Style = "Normal"
$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0842"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -2.31%  "
$c.Style = "Normal"

$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = "Chainlink"
$c.Style = "Normal"
$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "18.43"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -4.79%  "
$c.Style = "Normal"

$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = "WrappedliquidstakedEther2.0"
$c.Style = "Normal"
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.382.11"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -2.37%  "
$c.Style = "Normal"

$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "Polkadot"
$c.Style = "Normal"
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.41"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -4.79%  "
$c.Style = "Normal"

$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "WrappedEther"
$c.Style = "Normal"
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.922.32"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -1.98%  "
$c.Style = "Normal"

$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = "Polygon"
$c.Style = "Normal"
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.938"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -8.91%  "
$c.Style = "Normal"

$c = $ws.Range("B18")
$c.NumberFormat = "@"
$c.Value = "WrappedBTC"
$c.Style = "Normal"
$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "51.322.79"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -1.89%  "
$c.Style = "Normal"

$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "ImmutableX"
$c.Style = "Normal"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.32"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -6.04%  "
$c.Style = "Normal"

$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.Style = "Normal"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.34"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -2.97%  "
$c.Style = "Normal"

$c = $ws.Range("B21")
$c.NumberFormat = "@"
$c.Value = "InternetComputer(DFINITY)"
$c.Style = "Normal"
$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.08"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -4.25%  "
$c.Style = "Normal"

$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.Style = "Normal"
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.0₃0948"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -2.91%  "
$c.Style = "Normal"

$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = "Litecoin"
$c.Style = "Normal"
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "68.70"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -1.31%  "
$c.Style = "Normal"

$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.Style = "Normal"
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "261.64"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -1.13%  "
$c.Style = "Normal"

$c = $ws.Range("B25")
$c.NumberFormat = "@"
$c.Value = "PancakeSwap"
$c.Style = "Normal"
$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -1.52%  "
$c.Style = "Normal"

$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = "Kaspa"
$c.Style = "Normal"
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.172"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -4.45%  "
$c.Style = "Normal"

$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = "LEO"
$c.Style = "Normal"
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -3.89%  "
$c.Style = "Normal"

$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = "Dai"
$c.Style = "Normal"
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.Style = "Normal"

$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = "EthereumClassic"
$c.Style = "Normal"
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "25.92"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -3.89%  "
$c.Style = "Normal"

$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c.Style = "Normal"
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +7.85%  "
$c.Style = "Normal"

$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c.Style = "Normal"
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.18"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -5.33%  "
$c.Style = "Normal"

$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = "Hedera"
$c.Style = "Normal"
$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -4.46%  "
$c.Style = "Normal"

$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = "Cosmos"
$c.Style = "Normal"
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "9.93"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -4.24%  "
$c.Style = "Normal"

$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = "Toncoin"
$c.Style = "Normal"
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.12"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -3.52%  "
$c.Style = "Normal"

$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = "InjectiveProtocol"
$c.Style = "Normal"
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "34.80"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -4.98%  "
$c.Style = "Normal"

$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = "OKB"
$c.Style = "Normal"
$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "51.09"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +0.26%  "
$c.Style = "Normal"

$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "FirstDigitalUSD"
$c.Style = "Normal"
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +0.54%  "
$c.Style = "Normal"

$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.Style = "Normal"
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0427"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -4.29%  "
$c.Style = "Normal"

$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = "LidoDAOToken"
$c.Style = "Normal"
$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -5.80%  "
$c.Style = "Normal"

$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = "Celestia"
$c.Style = "Normal"
$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "17.21"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -4.30%  "
$c.Style = "Normal"

$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = "Stacks"
$c.Style = "Normal"
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -4.88%  "
$c.Style = "Normal"

$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = "ARBITRUM"
$c.Style = "Normal"
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -6.25%  "
$c.Style = "Normal"

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c.Style = "Normal"
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -3.84%  "
$c.Style = "Normal"

$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = "EnergySwap"
$c.Style = "Normal"
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "22.03"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -3.66%  "
$c.Style = "Normal"

$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.Style = "Normal"
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "119.57"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -3.84%  "
$c.Style = "Normal"

$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = "WEMIXToken"
$c.Style = "Normal"
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.09"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -2.70%  "
$c.Style = "Normal"

$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = "Maker"
$c.Style = "Normal"
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.032.82"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -4.39%  "
$c.Style = "Normal"

$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = "ApeXProtocol"
$c.Style = "Normal"
$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -3.69%  "
$c.Style = "Normal"

$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c.Style = "Normal"
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -5.43%  "
$c.Style = "Normal"

$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = "TheGraph"
$c.Style = "Normal"
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.260"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +4.87%  "
$c.Style = "Normal"

$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = "RocketPoolETH"
$c.Style = "Normal"
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.215.55"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -2.14%  "
$c.Style = "Normal"
